# Update "想去人数" (want-to-go count, column F) values across the four
# worksheets of the workbook to reflect the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1049
$ws.Range("F3").Value = 653
$ws.Range("F4").Value = 1453
$ws.Range("F5").Value = 61
$ws.Range("F6").Value = 3197
$ws.Range("F7").Value = 8
$ws.Range("F8").Value = 585
$ws.Range("F9").Value = 2133
$ws.Range("F10").Value = 456
$ws.Range("F11").Value = 386
$ws.Range("F12").Value = 228
$ws.Range("F14").Value = 256
$ws.Range("F16").Value = 1044
$ws.Range("F17").Value = 414
$ws.Range("F18").Value = 69
$ws.Range("F19").Value = 170
$ws.Range("F20").Value = 4185
$ws.Range("F21").Value = 1236
$ws.Range("F22").Value = 3254
$ws.Range("F23").Value = 316
$ws.Range("F24").Value = 119
$ws.Range("F25").Value = 3070
$ws.Range("F26").Value = 4718
$ws.Range("F27").Value = 120
$ws.Range("F28").Value = 959
$ws.Range("F29").Value = 517
$ws.Range("F30").Value = 3080
$ws.Range("F31").Value = 309
$ws.Range("F33").Value = 119
$ws.Range("F34").Value = 82
$ws.Range("F35").Value = 558
$ws.Range("F36").Value = 1110
$ws.Range("F37").Value = 1360
$ws.Range("F38").Value = 106
$ws.Range("F39").Value = 1259
$ws.Range("F40").Value = 805
$ws.Range("F42").Value = 752
$ws.Range("F43").Value = 480
$ws.Range("F45").Value = 235
$ws.Range("F46").Value = 50
$ws.Range("F47").Value = 101
$ws.Range("F48").Value = 353
$ws.Range("F49").Value = 3682

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 977
$ws.Range("F8").Value = 34
$ws.Range("F11").Value = 6
$ws.Range("F22").Value = 28
$ws.Range("F25").Value = 45

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1907

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1908
$ws.Range("F4").Value = 653
$ws.Range("F5").Value = 1453
$ws.Range("F6").Value = 3197
$ws.Range("F8").Value = 2133
$ws.Range("F9").Value = 456
$ws.Range("F10").Value = 386
$ws.Range("F12").Value = 228
$ws.Range("F13").Value = 977
$ws.Range("F14").Value = 34
$ws.Range("F16").Value = 256
$ws.Range("F17").Value = 1044
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = 414
$ws.Range("F20").Value = 170
$ws.Range("F21").Value = 4185
$ws.Range("F23").Value = 1236
$ws.Range("F25").Value = 3254
$ws.Range("F26").Value = 3070
$ws.Range("F27").Value = 4718
$ws.Range("F28").Value = 959
$ws.Range("F29").Value = 3080
$ws.Range("F30").Value = 309
$ws.Range("F32").Value = 119
$ws.Range("F33").Value = 82
$ws.Range("F34").Value = 558
$ws.Range("F35").Value = 1110
$ws.Range("F36").Value = 1360
$ws.Range("F37").Value = 106
$ws.Range("F38").Value = 1259
$ws.Range("F39").Value = 805
$ws.Range("F41").Value = 480
$ws.Range("F44").Value = 28
$ws.Range("F45").Value = 235
$ws.Range("F46").Value = 50
$ws.Range("F47").Value = 101
$ws.Range("F48").Value = 353
$ws.Range("F49").Value = 3682
